# RSTK-9576-SYDATA: Delete WO picklist transaction - template fix-up.
#
# 1. Rename the "Pro-Lot Track (Lot Track)" list value used in the
#    "Create WO" sheet (cell B2) to "Pro-SYDATA1 (Lot track)".
# 2. Switch the sheet to Portrait page orientation for printing.
# 3. Leave the final active selection on B2 (the cell that was edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create WO")

$ws.Activate()

# 1. Update the picklist label text in B2.
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# 2. Set the page to print in portrait orientation.
$ws.PageSetup.Orientation = 1

# 3. Select B2 so it becomes the sheet's active cell/selection.
$ws.Range("B2").Select()
